# Bulgaria First League - update of 28-05-2024 19:13
# 1) Rows 10 and 11 had their match (away team + stats) mixed up - swap
#    everything except the row's own id (col A) and home team (col E).
# 2) Two new finished... (pending stats) matches appended as rows 294/295.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Swap the data of rows 10 and 11 (columns B..AD), keep column A as-is.
# ---------------------------------------------------------------------
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

foreach ($col in $cols) {
    $addr10 = "$col" + "10"
    $addr11 = "$col" + "11"
    $v10 = $ws.Range($addr10).Value2
    $v11 = $ws.Range($addr11).Value2
    $ws.Range($addr10).Value = $v11
    $ws.Range($addr11).Value = $v10
}

# Column E (HomeTeam) must stay as it was before the swap.
$ws.Range("E10").Value = "Levski Sofia"
$ws.Range("E11").Value = "Cherno More Varna"

# ---------------------------------------------------------------------
# 2) Append the two new matches (rows 294 and 295).
# ---------------------------------------------------------------------
$ws.Range("A294").Value = 292
$ws.Range("B294").Value = 8129620
$ws.Range("C294").Value = "Bulgaria First League"
$ws.Range("D294").Value2 = 45439.58333333334
$ws.Range("E294").Value = "Slavia Sofia"
$ws.Range("F294").Value = "CSKA 1948 Sofia"
$ws.Range("G294").Value = 0
$ws.Range("H294").Value = 2
$ws.Range("K294").Value = "A"
$ws.Range("L294").Value = 3.75
$ws.Range("M294").Value = 3.6
$ws.Range("N294").Value = 1.9
$ws.Range("O294").Value = 9.5
$ws.Range("P294").Value = 5.25
$ws.Range("Q294").Value = 1.285
$ws.Range("R294").Value = 1.5
$ws.Range("S294").Value = 1.925
$ws.Range("T294").Value = 1.925
$ws.Range("U294").Value = 2.75
$ws.Range("V294").Value = 1.875
$ws.Range("W294").Value = 1.975
$ws.Range("X294").Value = -1
$ws.Range("Y294").Value = -1
$ws.Range("Z294").Value = 0.2849999999999999
$ws.Range("AA294").Value = -1
$ws.Range("AB294").Value = 0.925
$ws.Range("AC294").Value = -1
$ws.Range("AD294").Value = 0.9750000000000001

$ws.Range("A295").Value = 293
$ws.Range("B295").Value = 8129621
$ws.Range("C295").Value = "Bulgaria First League"
$ws.Range("D295").Value2 = 45439.58333333334
$ws.Range("E295").Value = "Botev Plovdiv"
$ws.Range("F295").Value = "Arda Kardzhali"
$ws.Range("G295").Value = 0
$ws.Range("H295").Value = 1
$ws.Range("K295").Value = "A"
$ws.Range("L295").Value = 4.333
$ws.Range("M295").Value = 4
$ws.Range("N295").Value = 1.727
$ws.Range("O295").Value = 8.5
$ws.Range("P295").Value = 5.25
$ws.Range("Q295").Value = 1.333
$ws.Range("R295").Value = 1.5
$ws.Range("S295").Value = 1.875
$ws.Range("T295").Value = 1.975
$ws.Range("U295").Value = 2.75
$ws.Range("V295").Value = 1.85
$ws.Range("W295").Value = 2
$ws.Range("X295").Value = -1
$ws.Range("Y295").Value = -1
$ws.Range("Z295").Value = 0.333
$ws.Range("AA295").Value = 0.875
$ws.Range("AB295").Value = -1
$ws.Range("AC295").Value = -1
$ws.Range("AD295").Value = 1

# Copy the formatting (bold/border/center for id, custom date format for
# Date) from the last pre-existing row down onto the two new rows.
$ws.Range("A293").Copy() | Out-Null
$ws.Range("A294:A295").PasteSpecial(-4122) | Out-Null

$ws.Range("D293").Copy() | Out-Null
$ws.Range("D294:D295").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
